$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: propagate row banding (border/number-format) styles into the new rows by
# copying formats from existing rows that already carry the matching style pattern for
# that banding parity / column "shape" (J:L filled vs M:O filled vs both). Done before any
# values are written and before row 517 is touched, since row 517s current formatting is
# itself used as a template below.

$ws.Range("A4:L4").Copy()
$ws.Range("A518:L518").PasteSpecial(-4122)
$ws.Range("A3:I3").Copy()
$ws.Range("A519:I519").PasteSpecial(-4122)
$ws.Range("M3:O3").Copy()
$ws.Range("M519:O519").PasteSpecial(-4122)
$ws.Range("A4:L4").Copy()
$ws.Range("A520:L520").PasteSpecial(-4122)
$ws.Range("A3:I3").Copy()
$ws.Range("A521:I521").PasteSpecial(-4122)
$ws.Range("M3:O3").Copy()
$ws.Range("M521:O521").PasteSpecial(-4122)
$ws.Range("A6:I6").Copy()
$ws.Range("A522:I522").PasteSpecial(-4122)
$ws.Range("M6:O6").Copy()
$ws.Range("M522:O522").PasteSpecial(-4122)
$ws.Range("A3:I3").Copy()
$ws.Range("A523:I523").PasteSpecial(-4122)
$ws.Range("M3:O3").Copy()
$ws.Range("M523:O523").PasteSpecial(-4122)
$ws.Range("A6:I6").Copy()
$ws.Range("A524:I524").PasteSpecial(-4122)
$ws.Range("M6:O6").Copy()
$ws.Range("M524:O524").PasteSpecial(-4122)
$ws.Range("A5:L5").Copy()
$ws.Range("A525:L525").PasteSpecial(-4122)
$ws.Range("A6:I6").Copy()
$ws.Range("A526:I526").PasteSpecial(-4122)
$ws.Range("M6:O6").Copy()
$ws.Range("M526:O526").PasteSpecial(-4122)
$ws.Range("A5:L5").Copy()
$ws.Range("A527:L527").PasteSpecial(-4122)
$ws.Range("A4:L4").Copy()
$ws.Range("A528:L528").PasteSpecial(-4122)
$ws.Range("A3:I3").Copy()
$ws.Range("A529:I529").PasteSpecial(-4122)
$ws.Range("M3:O3").Copy()
$ws.Range("M529:O529").PasteSpecial(-4122)
$ws.Range("A6:I6").Copy()
$ws.Range("A530:I530").PasteSpecial(-4122)
$ws.Range("M6:O6").Copy()
$ws.Range("M530:O530").PasteSpecial(-4122)

# Row 531 becomes the new last row of the table: its A:L banding matches the existing
# "last row of a JKL-shaped stretch" style (e.g. row 71), and its M:O filler cells reuse
# the same blank "no-border filler" style that row 517s J:K:L currently have (grabbed
# here, before that formatting is cleared off row 517 in step 2).
$ws.Range("A71:L71").Copy()
$ws.Range("A531:L531").PasteSpecial(-4122)
$ws.Range("J517:L517").Copy()
$ws.Range("M531:O531").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Step 2: row 517 is no longer the last row of the table, so (matching every other
# non-last data row already in the sheet) Excel drops its placeholder empty J:K:L filler cells.
$ws.Range("J517:L517").Clear()

# --- Step 3: write the 14 new Google-Form response rows (518-531).

# Row 518
$ws.Range("A518").Value = "45599.283865949074"
$ws.Range("B518").Value = 'hlu20242513@gmail.com'
$ws.Range("C518").Value = '미디어스쿨'
$ws.Range("D518").Value = 20242513
$ws.Range("E518").Value = '김예준'
$ws.Range("F518").Value = '가짜약 대조군'
$ws.Range("G518").Value = '28 vs 25'
$ws.Range("H518").Value = '플라시보 컨트롤 설계의 생리식염수 접종 집단'
$ws.Range("I518").Value = 'Red'
$ws.Range("J518").Value = '나. 5센트'
$ws.Range("K518").Value = '나. 5분'
$ws.Range("L518").Value = '가. 24일'

# Row 519
$ws.Range("A519").Value = "45599.289407013886"
$ws.Range("B519").Value = 'cbh3trust4@naver.com'
$ws.Range("C519").Value = '법학과'
$ws.Range("D519").Value = 20242749
$ws.Range("E519").Value = '조정현'
$ws.Range("F519").Value = '가짜약 대조군'
$ws.Range("G519").Value = '28 vs 25'
$ws.Range("H519").Value = 'NFIP 설계의 백신 접종 집단'
$ws.Range("I519").Value = 'Black'
$ws.Range("M519").Value = '가. 5센트'
$ws.Range("N519").Value = '가. 5분'
$ws.Range("O519").Value = '가. 47일'

# Row 520
$ws.Range("A520").Value = "45599.295889583329"
$ws.Range("B520").Value = 'hyeruys2005@naver.com'
$ws.Range("C520").Value = '미디어스쿨'
$ws.Range("D520").Value = 20242575
$ws.Range("E520").Value = '정윤수'
$ws.Range("F520").Value = '가짜약 대조군'
$ws.Range("G520").Value = '28 vs 25'
$ws.Range("H520").Value = 'NFIP 설계의 대조군 집단'
$ws.Range("I520").Value = 'Red'
$ws.Range("J520").Value = '나. 5센트'
$ws.Range("K520").Value = '가. 100분'
$ws.Range("L520").Value = '나. 47일'

# Row 521
$ws.Range("A521").Value = "45599.298872222222"
$ws.Range("B521").Value = 'yeon4262@naver.com'
$ws.Range("C521").Value = '반도체디스플레이스쿨'
$ws.Range("D521").Value = 20223325
$ws.Range("E521").Value = '신수연'
$ws.Range("F521").Value = '가짜약 대조군'
$ws.Range("G521").Value = '25 vs 54'
$ws.Range("H521").Value = '플라시보 컨트롤 설계의 생리식염수 접종 집단'
$ws.Range("I521").Value = 'Black'
$ws.Range("M521").Value = '가. 5센트'
$ws.Range("N521").Value = '가. 5분'
$ws.Range("O521").Value = '가. 47일'

# Row 522
$ws.Range("A522").Value = "45599.332078356485"
$ws.Range("B522").Value = 'sjhaa303028@naver.com'
$ws.Range("C522").Value = '인공지능융합학부'
$ws.Range("D522").Value = 20246741
$ws.Range("E522").Value = '신중현'
$ws.Range("F522").Value = '이중눈가림'
$ws.Range("G522").Value = '28 vs 46'
$ws.Range("H522").Value = '플라시보 컨트롤 설계의 생리식염수 접종 집단'
$ws.Range("I522").Value = 'Black'
$ws.Range("M522").Value = '나. 10센트'
$ws.Range("N522").Value = '나. 100분'
$ws.Range("O522").Value = '가. 47일'

# Row 523
$ws.Range("A523").Value = "45599.338345543976"
$ws.Range("B523").Value = 'ian5791@naver.com'
$ws.Range("C523").Value = '중국학과'
$ws.Range("D523").Value = 20241520
$ws.Range("E523").Value = '박수현'
$ws.Range("F523").Value = '랜덤화'
$ws.Range("G523").Value = '28 vs 71'
$ws.Range("H523").Value = '플라시보 컨트롤 설계의 생리식염수 접종 집단'
$ws.Range("I523").Value = 'Black'
$ws.Range("M523").Value = '나. 10센트'
$ws.Range("N523").Value = '나. 100분'
$ws.Range("O523").Value = '나. 24일'

# Row 524
$ws.Range("A524").Value = "45599.35747077546"
$ws.Range("B524").Value = 'ac3512@naver.com'
$ws.Range("C524").Value = '융합신소재공학'
$ws.Range("D524").Value = 20216609
$ws.Range("E524").Value = '김태근'
$ws.Range("F524").Value = '랜덤화'
$ws.Range("G524").Value = '28 vs 71'
$ws.Range("H524").Value = 'NFIP 설계의 대조군 집단'
$ws.Range("I524").Value = 'Black'
$ws.Range("M524").Value = '나. 10센트'
$ws.Range("N524").Value = '나. 100분'
$ws.Range("O524").Value = '나. 24일'

# Row 525
$ws.Range("A525").Value = "45599.359303356483"
$ws.Range("B525").Value = 'dldpwls5245@naver.com'
$ws.Range("C525").Value = '법학과'
$ws.Range("D525").Value = 20182747
$ws.Range("E525").Value = '이예진'
$ws.Range("F525").Value = '가짜약 대조군'
$ws.Range("G525").Value = '28 vs 71'
$ws.Range("H525").Value = '플라시보 컨트롤 설계의 생리식염수 접종 집단'
$ws.Range("I525").Value = 'Red'
$ws.Range("J525").Value = '가. 10센트'
$ws.Range("K525").Value = '나. 5분'
$ws.Range("L525").Value = '가. 24일'

# Row 526
$ws.Range("A526").Value = "45599.392086342588"
$ws.Range("B526").Value = 'rhkddyd234@naver.com'
$ws.Range("C526").Value = '화학과'
$ws.Range("D526").Value = 20193420
$ws.Range("E526").Value = '이광용'
$ws.Range("F526").Value = '이중눈가림'
$ws.Range("G526").Value = '28 vs 46'
$ws.Range("H526").Value = 'NFIP 설계의 대조군 집단'
$ws.Range("I526").Value = 'Black'
$ws.Range("M526").Value = '가. 5센트'
$ws.Range("N526").Value = '나. 100분'
$ws.Range("O526").Value = '나. 24일'

# Row 527
$ws.Range("A527").Value = "45599.413209131948"
$ws.Range("B527").Value = 'jeongminyoung@naver.com'
$ws.Range("C527").Value = '간호학과'
$ws.Range("D527").Value = 20246281
$ws.Range("E527").Value = '정민영'
$ws.Range("F527").Value = '랜덤화'
$ws.Range("G527").Value = '28 vs 71'
$ws.Range("H527").Value = 'NFIP 설계의 대조군 집단'
$ws.Range("I527").Value = 'Red'
$ws.Range("J527").Value = '나. 5센트'
$ws.Range("K527").Value = '나. 5분'
$ws.Range("L527").Value = '가. 24일'

# Row 528
$ws.Range("A528").Value = "45599.421258125003"
$ws.Range("B528").Value = 'jhkm7400@gmail.com'
$ws.Range("C528").Value = '경영대학'
$ws.Range("D528").Value = 20242932
$ws.Range("E528").Value = '김유건'
$ws.Range("F528").Value = '랜덤화'
$ws.Range("G528").Value = '28 vs 71'
$ws.Range("H528").Value = 'NFIP 설계의 대조군 집단'
$ws.Range("I528").Value = 'Red'
$ws.Range("J528").Value = '가. 10센트'
$ws.Range("K528").Value = '나. 5분'
$ws.Range("L528").Value = '가. 24일'

# Row 529
$ws.Range("A529").Value = "45599.422794745369"
$ws.Range("B529").Value = 'ghskfen1215@naver.com'
$ws.Range("C529").Value = '인공지능융합학부'
$ws.Range("D529").Value = 20236705
$ws.Range("E529").Value = '기정윤'
$ws.Range("F529").Value = '랜덤화'
$ws.Range("G529").Value = '28 vs 71'
$ws.Range("H529").Value = 'NFIP 설계의 대조군 집단'
$ws.Range("I529").Value = 'Black'
$ws.Range("M529").Value = '가. 5센트'
$ws.Range("N529").Value = '가. 5분'
$ws.Range("O529").Value = '가. 47일'

# Row 530
$ws.Range("A530").Value = "45599.427128726849"
$ws.Range("B530").Value = 'a22234781@gmail.com'
$ws.Range("C530").Value = '일본학과'
$ws.Range("D530").Value = 20201076
$ws.Range("E530").Value = '이수민'
$ws.Range("F530").Value = '랜덤화'
$ws.Range("G530").Value = '28 vs 25'
$ws.Range("H530").Value = 'NFIP 설계의 대조군 집단'
$ws.Range("I530").Value = 'Black'
$ws.Range("M530").Value = '나. 10센트'
$ws.Range("N530").Value = '나. 100분'
$ws.Range("O530").Value = '나. 24일'

# Row 531
$ws.Range("A531").Value = "45599.435293645831"
$ws.Range("B531").Value = 'hg2635394@gmail.com'
$ws.Range("C531").Value = '철학과'
$ws.Range("D531").Value = 20230133
$ws.Range("E531").Value = '김현준'
$ws.Range("F531").Value = '랜덤화'
$ws.Range("G531").Value = '28 vs 71'
$ws.Range("H531").Value = 'NFIP 설계의 대조군 집단'
$ws.Range("I531").Value = 'Red'
$ws.Range("J531").Value = '가. 10센트'
$ws.Range("K531").Value = '가. 100분'
$ws.Range("L531").Value = '가. 24일'

# --- Step 4: grow the Excel Table ("Form_Responses1") to cover the new rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:O531"))

# --- Step 5: restore the active selection to where the user ended up after data entry.
$ws.Range("G535").Select()

Write-Host "edit complete"
